$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04936366666666667
$ws.Range("H2").Value = 0.148091
$ws.Range("I2").Value = 0.04616170608573571
$ws.Range("J2").Value = 0.0461617060857357
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9949870000000001
$ws.Range("N2").Value = 2.984961
$ws.Range("O2").Value = 0.03855738270564991
$ws.Range("P2").Value = 0.03855738270564991
$ws.Range("Q2").Value = 0.04911620660566667
$ws.Range("R2").Value = 0.4420458594510001
$ws.Range("S2").Value = 0.00177987456789344
$ws.Range("T2").Value = 0.00177987456789344
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04936366666666667
$ws.Range("H3").Value = 0.148091
$ws.Range("I3").Value = 0.04616170608573571
$ws.Range("J3").Value = 0.0461617060857357
$ws.Range("O3").Value = 0.3282159160005915
$ws.Range("P3").Value = 0.3282159160005916
$ws.Range("Q3").Value = 0.4180968626584444
$ws.Range("R3").Value = 3.762871763926
$ws.Range("S3").Value = 0.01515100664707983
$ws.Range("T3").Value = 0.01515100664707983
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04936366666666667
$ws.Range("H4").Value = 0.148091
$ws.Range("I4").Value = 0.04616170608573571
$ws.Range("J4").Value = 0.0461617060857357
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04495399999999999
$ws.Range("N4").Value = 0.134862
$ws.Range("O4").Value = 0.001742041435867791
$ws.Range("P4").Value = 0.001742041435867791
$ws.Range("Q4").Value = 0.002219094271333333
$ws.Range("R4").Value = 0.019971848442
$ws.Range("S4").Value = 0.00008041560475170197
$ws.Range("T4").Value = 0.00008041560475170197
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04936366666666667
$ws.Range("H5").Value = 0.148091
$ws.Range("I5").Value = 0.04616170608573571
$ws.Range("J5").Value = 0.0461617060857357
$ws.Range("M5").Value = 15.96019966666667
$ws.Range("N5").Value = 47.880599
$ws.Range("O5").Value = 0.6184839868322428
$ws.Range("P5").Value = 0.6184839868322429
$ws.Range("Q5").Value = 0.7878539762787778
$ws.Range("R5").Value = 7.090685786509001
$ws.Range("S5").Value = 0.02855027601888403
$ws.Range("T5").Value = 0.02855027601888403
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04936366666666667
$ws.Range("H6").Value = 0.148091
$ws.Range("I6").Value = 0.04616170608573571
$ws.Range("J6").Value = 0.0461617060857357
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3268106666666666
$ws.Range("N6").Value = 0.980432
$ws.Range("O6").Value = 0.01266445083901121
$ws.Range("P6").Value = 0.01266445083901121
$ws.Range("Q6").Value = 0.01613257281244444
$ws.Range("R6").Value = 0.145193155312
$ws.Range("S6").Value = 0.0005846126573676846
$ws.Range("T6").Value = 0.0005846126573676846
$ws.Range("A7").Value = "ECs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 0.04936366666666667
$ws.Range("H7").Value = 0.148091
$ws.Range("I7").Value = 0.04616170608573571
$ws.Range("J7").Value = 0.0461617060857357
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.008676333333333333
$ws.Range("N7").Value = 0.026029
$ws.Range("O7").Value = 0.0003362221866367304
$ws.Range("P7").Value = 0.0003362221866367304
$ws.Range("Q7").Value = 0.0004282956265555555
$ws.Range("R7").Value = 0.003854660639
$ws.Range("S7").Value = 0.00001552058975902812
$ws.Range("T7").Value = 0.00001552058975902812
$ws.Range("D8").Value = "ECs"
$ws.Range("I8").Value = 0.008057748967298944
$ws.Range("J8").Value = 0.008057748967298944
$ws.Range("M8").Value = 0.9949870000000001
$ws.Range("N8").Value = 2.984961
$ws.Range("O8").Value = 0.03855738270564991
$ws.Range("P8").Value = 0.03855738270564991
$ws.Range("Q8").Value = 0.008573471316666668
$ws.Range("R8").Value = 0.07716124185000001
$ws.Range("S8").Value = 0.0003106857106782007
$ws.Range("T8").Value = 0.0003106857106782007
$ws.Range("D9").Value = "FAPs"
$ws.Range("I9").Value = 0.008057748967298944
$ws.Range("J9").Value = 0.008057748967298944
$ws.Range("M9").Value = 8.469728666666667
$ws.Range("N9").Value = 25.409186
$ws.Range("O9").Value = 0.3282159160005915
$ws.Range("P9").Value = 0.3282159160005916
$ws.Range("Q9").Value = 0.07298082867777778
$ws.Range("R9").Value = 0.6568274581000001
$ws.Range("S9").Value = 0.002644681458204843
$ws.Range("T9").Value = 0.002644681458204844
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("I10").Value = 0.008057748967298944
$ws.Range("J10").Value = 0.008057748967298944
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.04495399999999999
$ws.Range("N10").Value = 0.134862
$ws.Range("O10").Value = 0.001742041435867791
$ws.Range("P10").Value = 0.001742041435867791
$ws.Range("Q10").Value = 0.0003873536333333333
$ws.Range("R10").Value = 0.0034861827
$ws.Range("S10").Value = 0.00001403693258085566
$ws.Range("T10").Value = 0.00001403693258085566
$ws.Range("D11").Value = "MuSCs"
$ws.Range("I11").Value = 0.008057748967298944
$ws.Range("J11").Value = 0.008057748967298944
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 15.96019966666667
$ws.Range("N11").Value = 47.880599
$ws.Range("O11").Value = 0.6184839868322428
$ws.Range("P11").Value = 0.6184839868322429
$ws.Range("Q11").Value = 0.1375237204611111
$ws.Range("R11").Value = 1.23771348415
$ws.Range("S11").Value = 0.004983588706188438
$ws.Range("T11").Value = 0.004983588706188439
$ws.Range("A12").Value = "FAPs"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.008616666666666667
$ws.Range("H12").Value = 0.02585
$ws.Range("I12").Value = 0.008057748967298944
$ws.Range("J12").Value = 0.008057748967298944
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.3268106666666666
$ws.Range("N12").Value = 0.980432
$ws.Range("O12").Value = 0.01266445083901121
$ws.Range("P12").Value = 0.01266445083901121
$ws.Range("Q12").Value = 0.002816018577777777
$ws.Range("R12").Value = 0.0253441672
$ws.Range("S12").Value = 0.0001020469656694508
$ws.Range("T12").Value = 0.0001020469656694508
$ws.Range("A13").Value = "FAPs"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.008616666666666667
$ws.Range("H13").Value = 0.02585
$ws.Range("I13").Value = 0.008057748967298944
$ws.Range("J13").Value = 0.008057748967298944
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.008676333333333333
$ws.Range("N13").Value = 0.026029
$ws.Range("O13").Value = 0.0003362221866367304
$ws.Range("P13").Value = 0.0003362221866367304
$ws.Range("Q13").Value = 0.00007476107222222222
$ws.Range("R13").Value = 0.00067284965
$ws.Range("S13").Value = 0.000002709193977155107
$ws.Range("T13").Value = 0.000002709193977155108
$ws.Range("D14").Value = "ECs"
$ws.Range("G14").Value = 1.011383666666666
$ws.Range("H14").Value = 3.034151
$ws.Range("I14").Value = 0.9457805449469654
$ws.Range("J14").Value = 0.9457805449469653
$ws.Range("M14").Value = 0.9949870000000001
$ws.Range("N14").Value = 2.984961
$ws.Range("O14").Value = 0.03855738270564991
$ws.Range("P14").Value = 0.03855738270564991
$ws.Range("Q14").Value = 1.006313600345667
$ws.Range("R14").Value = 9.056822403110999
$ws.Range("S14").Value = 0.03646682242707827
$ws.Range("T14").Value = 0.03646682242707827
$ws.Range("D15").Value = "FAPs"
$ws.Range("G15").Value = 1.011383666666666
$ws.Range("H15").Value = 3.034151
$ws.Range("I15").Value = 0.9457805449469654
$ws.Range("J15").Value = 0.9457805449469653
$ws.Range("M15").Value = 8.469728666666667
$ws.Range("N15").Value = 25.409186
$ws.Range("O15").Value = 0.3282159160005915
$ws.Range("P15").Value = 0.3282159160005916
$ws.Range("Q15").Value = 8.566145234565109
$ws.Range("R15").Value = 77.095307111086
$ws.Range("S15").Value = 0.3104202278953069
$ws.Range("T15").Value = 0.3104202278953069
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 1.011383666666666
$ws.Range("H16").Value = 3.034151
$ws.Range("I16").Value = 0.9457805449469654
$ws.Range("J16").Value = 0.9457805449469653
$ws.Range("M16").Value = 0.04495399999999999
$ws.Range("N16").Value = 0.134862
$ws.Range("O16").Value = 0.001742041435867791
$ws.Range("P16").Value = 0.001742041435867791
$ws.Range("Q16").Value = 0.04546574135133332
$ws.Range("R16").Value = 0.4091916721619999
$ws.Range("S16").Value = 0.001647588898535233
$ws.Range("T16").Value = 0.001647588898535233
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Ptprz1"
$ws.Range("C17").Value = "Ncam1"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.011383666666666
$ws.Range("H17").Value = 3.034151
$ws.Range("I17").Value = 0.9457805449469654
$ws.Range("J17").Value = 0.9457805449469653
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 15.96019966666667
$ws.Range("N17").Value = 47.880599
$ws.Range("O17").Value = 0.6184839868322428
$ws.Range("P17").Value = 0.6184839868322429
$ws.Range("Q17").Value = 16.14188525960544
$ws.Range("R17").Value = 145.276967336449
$ws.Range("S17").Value = 0.5849501221071703
$ws.Range("T17").Value = 0.5849501221071705
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Ptprz1"
$ws.Range("C18").Value = "Ncam1"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.011383666666666
$ws.Range("H18").Value = 3.034151
$ws.Range("I18").Value = 0.9457805449469654
$ws.Range("J18").Value = 0.9457805449469653
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.3268106666666666
$ws.Range("N18").Value = 0.980432
$ws.Range("O18").Value = 0.01266445083901121
$ws.Range("P18").Value = 0.01266445083901121
$ws.Range("Q18").Value = 0.330530970359111
$ws.Range("R18").Value = 2.974778733232
$ws.Range("S18").Value = 0.01197779121597408
$ws.Range("T18").Value = 0.01197779121597408
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Ptprz1"
$ws.Range("C19").Value = "Ncam1"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.011383666666666
$ws.Range("H19").Value = 3.034151
$ws.Range("I19").Value = 0.9457805449469654
$ws.Range("J19").Value = 0.9457805449469653
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.008676333333333333
$ws.Range("N19").Value = 0.026029
$ws.Range("O19").Value = 0.0003362221866367304
$ws.Range("P19").Value = 0.0003362221866367304
$ws.Range("Q19").Value = 0.008775101819888886
$ws.Range("R19").Value = 0.07897591637899999
$ws.Range("S19").Value = 0.0003179924029005472
$ws.Range("T19").Value = 0.0003179924029005472
